# Add a new data row (row 8) to Sheet1 for patient Pt0 / Pt0_Bcells / chrY,
# mirroring the existing rows' structure (ref/alt depth counts + VAF formula),
# and move the active selection to H9 (matches the author's final cursor
# position recorded in the sheet view).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Pt0"
$ws.Range("B8").Value = "Pt0_Bcells"
$ws.Range("C8").Value = "chrY"
$ws.Range("D8").Value = 1043
$ws.Range("E8").Value = 41
$ws.Range("F8").Value = 9
$ws.Range("G8").Formula = "=(F8/(F8+E8))"

$ws.Range("H9").Select()
